# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the e815c2b6-a9ff-4e28-9054-ee6d26f5fb0f report row(s) now that the
# handback report has been regenerated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the e815c2b6... row.
$wsOverview.Range("G4").Value = "2016-08-16 00:42:02"

# zh-cn detail sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the e815c2b6... row.
$wsZhCn.Range("H4").Value = "2016-08-16 00:41:56"
$wsZhCn.Range("K4").Value = "2016-08-16 00:42:28"

# de-de detail sheet: Correspond Handoff Datetime (mirrors the Overview date)
# and Correspond Handback DateTime for the e815c2b6... row.
$wsDeDe.Range("H4").Value = "2016-08-16 00:42:02"
$wsDeDe.Range("K4").Value = "2016-08-16 00:42:35"
